$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @{
    4  = @{ C = 0.04321632163216321;  E = 0.03375937593759375;  G = 0.06677467746774676;  H = -0.06741074107410741;  I = -0.03313531353135313;  J = -0.3778668662055488 }
    5  = @{ C = -0.1412301230123012;  E = -0.106006600660066;   G = 0.01342934293429343;  H = -0.06672667266726671;  I = -0.1426222622262226;   J = 0.005223273418582731 }
    6  = @{ C = 0.05843384338433843;  E = -0.132985298529853;   G = -0.03218721872187218; H = -0.0144974497449745;   I = -0.04927692769276927;  J = -0.09372439863746897 }
    7  = @{ C = -0.004956495649564956; E = -0.05807380738073806; G = -0.0329072907290729;  H = 0.006540654065406539;  I = 0.007704770477047704;  J = -0.04798318952882156 }
    8  = @{ C = 0.0454005400540054;   E = 0.1906630663066306;   G = -0.9257125712571256;  H = 0.9999999999999999;    I = 0.8634503450345034;    J = 0.3602495884266467 }
    9  = @{ C = 0.9729972997299728;   E = 0.1639243924392439;   G = 0.07768376837683767;  H = -0.0007320732073207321; I = 0.09438943894389439;  J = 0.0352105128723333 }
    10 = @{ C = -0.08392439243924392; E = -0.2229222922292229;  G = -0.1218721872187219;  H = 0.05572157215721572;   I = -0.03294329432943294;  J = 0.0759568540743728 }
    11 = @{ C = -0.06389438943894389; E = 0.000588058805880588; G = 0.2182778277827783;   H = -0.2049324932493249;   I = -0.1291929192919292;   J = -0.2098084947630021 }
    12 = @{ C = -0.02811881188118811; E = -0.0852085208520852;  G = 0.05368136813681367;  H = -0.09008100810081007;  I = -0.1041584158415841;   J = 0.02466779529328369 }
    13 = @{ C = 0.2622622262226222;   E = 0.8611461146114611;   G = 0.3498829882988299;   H = -0.02807080708070807;  I = 0.4652265226522652;    J = -0.1093100694940685 }
    14 = @{ C = -0.1043744374437444;  E = -0.1631563156315631;  G = -0.08169216921692168; H = 0.03169516951695169;   I = -0.03866786678667866;  J = 0.06302188929095504 }
}

foreach ($row in $data.Keys) {
    $cols = $data[$row]
    foreach ($col in $cols.Keys) {
        $ws.Range("$col$row").Value = $cols[$col]
    }
}
